$wb = $excel.ActiveWorkbook

# Sheet "Overview": G3 holds the "Latest HO Xliff Generate Date" for the
# 33c72347... entry. Update the timestamp.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-03 02:49:17"

# Sheet "zh-cn": H3 = Correspond Handoff Datetime, K3 = Correspond Handback DateTime
# for the 33c72347... entry.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-03 02:49:13"
$wsZhCn.Range("K3").Value = "2016-09-03 02:49:31"

# Sheet "de-de": H3 = Correspond Handoff Datetime (shares the same text as
# Overview G3), K3 = Correspond Handback DateTime for the 33c72347... entry.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-03 02:49:17"
$wsDeDe.Range("K3").Value = "2016-09-03 02:49:40"
